$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.645.23'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.690.67'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.07'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.58'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('E9').Value = '  +5.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.08'
$ws.Range('E10').Value = '  +5.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.402'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '30.22'
$ws.Range('E13').Value = '  +4.35%  '
$ws.Range('E14').Value = '  +10.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.174.48'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.489.19'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.694.87'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.70'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '359.68'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.53'
$ws.Range('E21').Value = '  +3.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.13'
$ws.Range('E23').Value = '  +2.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.82'
$ws.Range('E24').Value = '  +3.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000107'
$ws.Range('E25').Value = '  +12.79%  '
$ws.Range('E26').Value = '  +2.95%  '
$ws.Range('E27').Value = '  -4.29%  '
$ws.Range('E28').Value = '  +3.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.26'
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('E30').Value = '  +5.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '534.94'
$ws.Range('E32').Value = '  +2.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.79'
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.67'
$ws.Range('E34').Value = '  +5.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.46'
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('E36').Value = '  +1.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.75'
$ws.Range('E37').Value = '  +2.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.86'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.00'
$ws.Range('E39').Value = '  -1.90%  '
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '168.80'
$ws.Range('E42').Value = '  +2.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.60'
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('E44').Value = '  +2.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0615'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.50'
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('E47').Value = '  +2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0267'
$ws.Range('E48').Value = '  +4.53%  '
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.05'
$ws.Range('E50').Value = '  +8.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0981'
$ws.Range('E51').Value = '  -0.26%  '
